$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3265.7778
$ws.Range("I76").Value = 3199
$ws.Range("J76").Value = 3499.5
$ws.Range("K76").Value = 3199
$ws.Range("L76").Value = 3499.5
$ws.Range("M76").Value = -2884
$ws.Range("N76").Value = -4129.5

# Hunk 1: ALC!row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3265.7778
$ws.Range("I79").Value = 3199
$ws.Range("J79").Value = 3499.5
$ws.Range("K79").Value = 3199
$ws.Range("L79").Value = 3499.5
$ws.Range("M79").Value = -2107
$ws.Range("N79").Value = -5683.5

# Hunk 2: ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 565.5714
$ws.Range("I98").Value = 643.1667
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 643.1667
$ws.Range("L98").Value = 100
$ws.Range("M98").Value = 854.8333
$ws.Range("N98").Value = -3096

# Hunk 3: ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 565.5714
$ws.Range("I122").Value = 643.1667
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 1929.5001
$ws.Range("L122").Value = 300
$ws.Range("M122").Value = 520.4999
$ws.Range("N122").Value = -5200

# Hunk 4: ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1942.2115
$ws.Range("I137").Value = 1989.4546
$ws.Range("K137").Value = 5968.3638
$ws.Range("M137").Value = -3418.3638

# Hunk 5: ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2955.2122
$ws.Range("J138").Value = 3532.95
$ws.Range("L138").Value = 10598.85
$ws.Range("N138").Value = -20878.85

# Hunk 6: ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2910.1296
$ws.Range("I32").Value = 1869.6222
$ws.Range("J32").Value = 8112.6665
$ws.Range("K32").Value = 1869.6222
$ws.Range("L32").Value = 8112.6665
$ws.Range("M32").Value = -1582.6222
$ws.Range("N32").Value = -8686.666499999999

# Hunk 7: ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3929.5
$ws.Range("I61").Value = 3929.4285
$ws.Range("J61").Value = 3929.5715
$ws.Range("K61").Value = 3929.4285
$ws.Range("L61").Value = 3929.5715
$ws.Range("M61").Value = -3717.4285
$ws.Range("N61").Value = -4353.5715

# Hunk 8: ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 53312
$ws.Range("I132").Value = 3303.1428
$ws.Range("K132").Value = 9909.428400000001
$ws.Range("M132").Value = -7379.428400000001

# Hunk 9: ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3929.5
$ws.Range("I136").Value = 3929.4285
$ws.Range("J136").Value = 3929.5715
$ws.Range("K136").Value = 11788.2855
$ws.Range("L136").Value = 11788.7145
$ws.Range("M136").Value = -9238.2855
$ws.Range("N136").Value = -16888.7145

# Hunk 10: BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1682.4736
$ws.Range("I105").Value = 1520.875
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 1520.875
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = 226.125
$ws.Range("N105").Value = -5294

# Hunk 11: BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1500
$ws.Range("J107").Value = 1500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340

# Hunk 12: BSM!row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3503.7036
$ws.Range("I134").Value = 3882.318
$ws.Range("J134").Value = 1837.8
$ws.Range("K134").Value = 11646.954
$ws.Range("L134").Value = 5513.4
$ws.Range("M134").Value = -9111.954000000002
$ws.Range("N134").Value = -10583.4

# Hunk 13: CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23715.25
$ws.Range("I31").Value = 152750.5
$ws.Range("J31").Value = 5281.643
$ws.Range("K31").Value = 152750.5
$ws.Range("L31").Value = 5281.643
$ws.Range("M31").Value = -152455.5
$ws.Range("N31").Value = -5871.643

# Hunk 14: CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23715.25
$ws.Range("I34").Value = 152750.5
$ws.Range("J34").Value = 5281.643
$ws.Range("K34").Value = 152750.5
$ws.Range("L34").Value = 5281.643
$ws.Range("M34").Value = -152548.5
$ws.Range("N34").Value = -5685.643

# Hunk 15: CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15734.703
$ws.Range("I132").Value = 23675.045
$ws.Range("J132").Value = 4088.8667
$ws.Range("K132").Value = 71025.13499999999
$ws.Range("L132").Value = 12266.6001
$ws.Range("M132").Value = -68495.13499999999
$ws.Range("N132").Value = -17326.6001

# Hunk 16: CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 967.5714
$ws.Range("I134").Value = 965.0769
$ws.Range("K134").Value = 2895.2307
$ws.Range("M134").Value = -360.2307000000001

# Hunk 17: CUL!row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6418
$ws.Range("J3").Value = 7985.2
$ws.Range("L3").Value = 23955.6
$ws.Range("N3").Value = -24179.6

# Hunk 18: CUL!row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1351.9032
$ws.Range("I68").Value = 586.55554
$ws.Range("J68").Value = 1665
$ws.Range("K68").Value = 1759.66662
$ws.Range("L68").Value = 4995
$ws.Range("M68").Value = -948.66662
$ws.Range("N68").Value = -6617

# Hunk 19: CUL!row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1351.9032
$ws.Range("I71").Value = 586.55554
$ws.Range("J71").Value = 1665
$ws.Range("K71").Value = 5278.99986
$ws.Range("L71").Value = 14985
$ws.Range("M71").Value = -1222.99986
$ws.Range("N71").Value = -23097

# Hunk 20: CUL!row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3057.805
$ws.Range("I107").Value = 8827.083000000001
$ws.Range("K107").Value = 26481.249
$ws.Range("M107").Value = -24561.249

# Hunk 21: CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 785.52
$ws.Range("J131").Value = 802.80414
$ws.Range("L131").Value = 2408.41242
$ws.Range("N131").Value = -12488.41242

# Hunk 22: CUL!row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8362.200000000001
$ws.Range("J137").Value = 3624.7222
$ws.Range("L137").Value = 10874.1666
$ws.Range("N137").Value = -21074.1666

# Hunk 23: GSM!row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 18000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 18000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -18576

# Hunk 24: GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9071.529
$ws.Range("J80").Value = 3862.625
$ws.Range("L80").Value = 3862.625
$ws.Range("N80").Value = -5858.625

# Hunk 25: GSM!row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 18000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 18000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 18000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -19996

# Hunk 26: GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 9071.529
$ws.Range("J83").Value = 3862.625
$ws.Range("L83").Value = 19313.125
$ws.Range("N83").Value = -29297.125

# Hunk 27: GSM!row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 18000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 18000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 54000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -63984

# Hunk 28: GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 31236.37
$ws.Range("I132").Value = 5263.7856
$ws.Range("J132").Value = 103959.6
$ws.Range("K132").Value = 15791.3568
$ws.Range("L132").Value = 311878.8
$ws.Range("M132").Value = -13261.3568
$ws.Range("N132").Value = -316938.8

# Hunk 29: LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5200.5
$ws.Range("I22").Value = 5200.5
$ws.Range("K22").Value = 5200.5
$ws.Range("M22").Value = -4905.5

# Hunk 30: LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5200.5
$ws.Range("I27").Value = 5200.5
$ws.Range("K27").Value = 5200.5
$ws.Range("M27").Value = -5093.5

# Hunk 31: LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 163185.58
$ws.Range("I40").Value = 225762.2
$ws.Range("K40").Value = 225762.2
$ws.Range("M40").Value = -225626.2

# Hunk 32: LTW!row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5357.65
$ws.Range("I61").Value = 2262.3333
$ws.Range("J61").Value = 10000.625
$ws.Range("K61").Value = 2262.3333
$ws.Range("L61").Value = 10000.625
$ws.Range("M61").Value = -2060.3333
$ws.Range("N61").Value = -10404.625

# Hunk 33: LTW!row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5357.65
$ws.Range("I113").Value = 2262.3333
$ws.Range("J113").Value = 10000.625
$ws.Range("K113").Value = 2262.3333
$ws.Range("L113").Value = 10000.625
$ws.Range("M113").Value = -92.33329999999978
$ws.Range("N113").Value = -14340.625

# Hunk 34: LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2466.1333
$ws.Range("I132").Value = 1927.5714
$ws.Range("J132").Value = 2937.375
$ws.Range("K132").Value = 5782.7142
$ws.Range("L132").Value = 8812.125
$ws.Range("M132").Value = -3252.7142
$ws.Range("N132").Value = -13872.125

# Hunk 35: LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 32400.059
$ws.Range("I136").Value = 57244.11
$ws.Range("J136").Value = 4450.5
$ws.Range("K136").Value = 171732.33
$ws.Range("L136").Value = 13351.5
$ws.Range("M136").Value = -169182.33
$ws.Range("N136").Value = -18451.5

# Hunk 36: WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1361.12
$ws.Range("I136").Value = 880.5714
$ws.Range("J136").Value = 1972.7273
$ws.Range("K136").Value = 2641.7142
$ws.Range("L136").Value = 5918.1819
$ws.Range("M136").Value = -91.71420000000035
$ws.Range("N136").Value = -11018.1819
